# Re-order rows 6-15 of the category dictionary sheet so that the
# K_CRIMOFF entry (previously the last row, 15) moves up to row 6, and
# the remaining entries shift down, resulting in an alphabetically
# sorted list of KNr codes from row 6 through row 15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("K_CRIMOFF", "Straftaten", "Criminal offences"),
    @("K_KREIS", "Kreis", "County"),
    @("K_LAENDER", "Bundesland", "Federal state"),
    @("K_PM", "Feinstaub", "Fine particulate matter"),
    @("K_SEA", "Meer", "Sea"),
    @("K_SERIES", "Zeitreihe", "Time series"),
    @("K_SEX", "Geschlecht", "Sex"),
    @("K_SUBINDEX", "Subindikatoren", "Sub index"),
    @("K_TYPEAREA", "Art der Fläche", "Type of area"),
    @("K_URBAN", "Verstädterungsgrad", "Degree of urbanisation")
)

$startRow = 6
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
